# Update "想去人数" (number of people interested) figures for several
# manga/anime convention listings, as published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 144   # 143 -> 144
$ws1.Range("F3").Value = 460   # 457 -> 460
$ws1.Range("F8").Value = 8     # 7   -> 8
$ws1.Range("F9").Value = 146   # 138 -> 146

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 80    # 79 -> 80
$ws2.Range("F3").Value = 31    # 30 -> 31

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 144   # 143 -> 144
$ws4.Range("F3").Value = 80    # 79  -> 80
$ws4.Range("F4").Value = 460   # 457 -> 460
$ws4.Range("F9").Value = 8     # 7   -> 8
$ws4.Range("F10").Value = 146  # 138 -> 146
$ws4.Range("F11").Value = 31   # 30  -> 31
